$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'ECs'
$ws.Cells.Item(2, 2).Value = 'Il16'
$ws.Cells.Item(2, 3).Value = 'Kcnj15'
$ws.Cells.Item(2, 4).Value = 'ECs'
$ws.Cells.Item(2, 5).Value = [double]"3"
$ws.Cells.Item(2, 6).Value = [double]"1"
$ws.Cells.Item(2, 7).Value = [double]"4.141153333333333"
$ws.Cells.Item(2, 8).Value = [double]"12.42346"
$ws.Cells.Item(2, 9).Value = [double]"0.2530231305454066"
$ws.Cells.Item(2, 10).Value = [double]"0.2530231305454066"
$ws.Cells.Item(2, 11).Value = [double]"3"
$ws.Cells.Item(2, 12).Value = [double]"1"
$ws.Cells.Item(2, 13).Value = [double]"2.359186333333333"
$ws.Cells.Item(2, 14).Value = [double]"7.077559"
$ws.Cells.Item(2, 15).Value = [double]"0.243271749000506"
$ws.Cells.Item(2, 16).Value = [double]"0.243271749000506"
$ws.Cells.Item(2, 17).Value = [double]"9.769752348237775"
$ws.Cells.Item(2, 18).Value = [double]"87.92777113413999"
$ws.Cells.Item(2, 19).Value = [double]"0.06155337950536444"
$ws.Cells.Item(2, 20).Value = [double]"0.06155337950536444"

# Row 3
$ws.Cells.Item(3, 1).Value = 'ECs'
$ws.Cells.Item(3, 2).Value = 'Il16'
$ws.Cells.Item(3, 3).Value = 'Kcnj15'
$ws.Cells.Item(3, 4).Value = 'FAPs'
$ws.Cells.Item(3, 5).Value = [double]"3"
$ws.Cells.Item(3, 6).Value = [double]"1"
$ws.Cells.Item(3, 7).Value = [double]"4.141153333333333"
$ws.Cells.Item(3, 8).Value = [double]"12.42346"
$ws.Cells.Item(3, 9).Value = [double]"0.2530231305454066"
$ws.Cells.Item(3, 10).Value = [double]"0.2530231305454066"
$ws.Cells.Item(3, 11).Value = [double]"3"
$ws.Cells.Item(3, 12).Value = [double]"1"
$ws.Cells.Item(3, 13).Value = [double]"7.327491999999999"
$ws.Cells.Item(3, 14).Value = [double]"21.982476"
$ws.Cells.Item(3, 15).Value = [double]"0.7555875385682617"
$ws.Cells.Item(3, 16).Value = [double]"0.7555875385682617"
$ws.Cells.Item(3, 17).Value = [double]"30.34426792077333"
$ws.Cells.Item(3, 18).Value = [double]"273.09841128696"
$ws.Cells.Item(3, 19).Value = [double]"0.1911811244096397"
$ws.Cells.Item(3, 20).Value = [double]"0.1911811244096397"

# Row 4
$ws.Cells.Item(4, 1).Value = 'ECs'
$ws.Cells.Item(4, 2).Value = 'Il16'
$ws.Cells.Item(4, 3).Value = 'Kcnj15'
$ws.Cells.Item(4, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(4, 5).Value = [double]"3"
$ws.Cells.Item(4, 6).Value = [double]"1"
$ws.Cells.Item(4, 7).Value = [double]"4.141153333333333"
$ws.Cells.Item(4, 8).Value = [double]"12.42346"
$ws.Cells.Item(4, 9).Value = [double]"0.2530231305454066"
$ws.Cells.Item(4, 10).Value = [double]"0.2530231305454066"
$ws.Cells.Item(4, 11).Value = [double]"1"
$ws.Cells.Item(4, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(4, 13).Value = [double]"0.01106233333333333"
$ws.Cells.Item(4, 14).Value = [double]"0.033187"
$ws.Cells.Item(4, 15).Value = [double]"0.001140712431232264"
$ws.Cells.Item(4, 16).Value = [double]"0.001140712431232264"
$ws.Cells.Item(4, 17).Value = [double]"0.04581081855777777"
$ws.Cells.Item(4, 18).Value = [double]"0.41229736702"
$ws.Cells.Item(4, 19).Value = [double]"0.0002886266304024494"
$ws.Cells.Item(4, 20).Value = [double]"0.0002886266304024494"

# Row 5
$ws.Cells.Item(5, 1).Value = 'FAPs'
$ws.Cells.Item(5, 2).Value = 'Il16'
$ws.Cells.Item(5, 3).Value = 'Kcnj15'
$ws.Cells.Item(5, 4).Value = 'ECs'
$ws.Cells.Item(5, 5).Value = [double]"3"
$ws.Cells.Item(5, 6).Value = [double]"1"
$ws.Cells.Item(5, 7).Value = [double]"5.865491666666667"
$ws.Cells.Item(5, 8).Value = [double]"17.596475"
$ws.Cells.Item(5, 9).Value = [double]"0.3583796455306321"
$ws.Cells.Item(5, 10).Value = [double]"0.358379645530632"
$ws.Cells.Item(5, 11).Value = [double]"3"
$ws.Cells.Item(5, 12).Value = [double]"1"
$ws.Cells.Item(5, 13).Value = [double]"2.359186333333333"
$ws.Cells.Item(5, 14).Value = [double]"7.077559"
$ws.Cells.Item(5, 15).Value = [double]"0.243271749000506"
$ws.Cells.Item(5, 16).Value = [double]"0.243271749000506"
$ws.Cells.Item(5, 17).Value = [double]"13.83778777828056"
$ws.Cells.Item(5, 18).Value = [double]"124.540090004525"
$ws.Cells.Item(5, 19).Value = [double]"0.08718364317441825"
$ws.Cells.Item(5, 20).Value = [double]"0.08718364317441823"

# Row 6
$ws.Cells.Item(6, 1).Value = 'FAPs'
$ws.Cells.Item(6, 2).Value = 'Il16'
$ws.Cells.Item(6, 3).Value = 'Kcnj15'
$ws.Cells.Item(6, 4).Value = 'FAPs'
$ws.Cells.Item(6, 5).Value = [double]"3"
$ws.Cells.Item(6, 6).Value = [double]"1"
$ws.Cells.Item(6, 7).Value = [double]"5.865491666666667"
$ws.Cells.Item(6, 8).Value = [double]"17.596475"
$ws.Cells.Item(6, 9).Value = [double]"0.3583796455306321"
$ws.Cells.Item(6, 10).Value = [double]"0.358379645530632"
$ws.Cells.Item(6, 11).Value = [double]"3"
$ws.Cells.Item(6, 12).Value = [double]"1"
$ws.Cells.Item(6, 13).Value = [double]"7.327491999999999"
$ws.Cells.Item(6, 14).Value = [double]"21.982476"
$ws.Cells.Item(6, 15).Value = [double]"0.7555875385682617"
$ws.Cells.Item(6, 16).Value = [double]"0.7555875385682617"
$ws.Cells.Item(6, 17).Value = [double]"42.97934326356667"
$ws.Cells.Item(6, 18).Value = [double]"386.8140893721"
$ws.Cells.Item(6, 19).Value = [double]"0.2707871942394564"
$ws.Cells.Item(6, 20).Value = [double]"0.2707871942394564"

# Row 7
$ws.Cells.Item(7, 1).Value = 'FAPs'
$ws.Cells.Item(7, 2).Value = 'Il16'
$ws.Cells.Item(7, 3).Value = 'Kcnj15'
$ws.Cells.Item(7, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(7, 5).Value = [double]"3"
$ws.Cells.Item(7, 6).Value = [double]"1"
$ws.Cells.Item(7, 7).Value = [double]"5.865491666666667"
$ws.Cells.Item(7, 8).Value = [double]"17.596475"
$ws.Cells.Item(7, 9).Value = [double]"0.3583796455306321"
$ws.Cells.Item(7, 10).Value = [double]"0.358379645530632"
$ws.Cells.Item(7, 11).Value = [double]"1"
$ws.Cells.Item(7, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(7, 13).Value = [double]"0.01106233333333333"
$ws.Cells.Item(7, 14).Value = [double]"0.033187"
$ws.Cells.Item(7, 15).Value = [double]"0.001140712431232264"
$ws.Cells.Item(7, 16).Value = [double]"0.001140712431232264"
$ws.Cells.Item(7, 17).Value = [double]"0.06488602398055557"
$ws.Cells.Item(7, 18).Value = [double]"0.5839742158250001"
$ws.Cells.Item(7, 19).Value = [double]"0.0004088081167574044"
$ws.Cells.Item(7, 20).Value = [double]"0.0004088081167574043"

# Row 8
$ws.Cells.Item(8, 1).Value = 'MuSCs'
$ws.Cells.Item(8, 2).Value = 'Il16'
$ws.Cells.Item(8, 3).Value = 'Kcnj15'
$ws.Cells.Item(8, 4).Value = 'ECs'
$ws.Cells.Item(8, 5).Value = [double]"3"
$ws.Cells.Item(8, 6).Value = [double]"1"
$ws.Cells.Item(8, 7).Value = [double]"0.467591"
$ws.Cells.Item(8, 8).Value = [double]"1.402773"
$ws.Cells.Item(8, 9).Value = [double]"0.02856965900840602"
$ws.Cells.Item(8, 10).Value = [double]"0.02856965900840601"
$ws.Cells.Item(8, 11).Value = [double]"3"
$ws.Cells.Item(8, 12).Value = [double]"1"
$ws.Cells.Item(8, 13).Value = [double]"2.359186333333333"
$ws.Cells.Item(8, 14).Value = [double]"7.077559"
$ws.Cells.Item(8, 15).Value = [double]"0.243271749000506"
$ws.Cells.Item(8, 16).Value = [double]"0.243271749000506"
$ws.Cells.Item(8, 17).Value = [double]"1.103134296789667"
$ws.Cells.Item(8, 18).Value = [double]"9.928208671107001"
$ws.Cells.Item(8, 19).Value = [double]"0.006950190915322995"
$ws.Cells.Item(8, 20).Value = [double]"0.006950190915322993"

# Row 9
$ws.Cells.Item(9, 1).Value = 'MuSCs'
$ws.Cells.Item(9, 2).Value = 'Il16'
$ws.Cells.Item(9, 3).Value = 'Kcnj15'
$ws.Cells.Item(9, 4).Value = 'FAPs'
$ws.Cells.Item(9, 5).Value = [double]"3"
$ws.Cells.Item(9, 6).Value = [double]"1"
$ws.Cells.Item(9, 7).Value = [double]"0.467591"
$ws.Cells.Item(9, 8).Value = [double]"1.402773"
$ws.Cells.Item(9, 9).Value = [double]"0.02856965900840602"
$ws.Cells.Item(9, 10).Value = [double]"0.02856965900840601"
$ws.Cells.Item(9, 11).Value = [double]"3"
$ws.Cells.Item(9, 12).Value = [double]"1"
$ws.Cells.Item(9, 13).Value = [double]"7.327491999999999"
$ws.Cells.Item(9, 14).Value = [double]"21.982476"
$ws.Cells.Item(9, 15).Value = [double]"0.7555875385682617"
$ws.Cells.Item(9, 16).Value = [double]"0.7555875385682617"
$ws.Cells.Item(9, 17).Value = [double]"3.426269311772"
$ws.Cells.Item(9, 18).Value = [double]"30.836423805948"
$ws.Cells.Item(9, 19).Value = [double]"0.02158687832789607"
$ws.Cells.Item(9, 20).Value = [double]"0.02158687832789606"

# Row 10
$ws.Cells.Item(10, 1).Value = 'MuSCs'
$ws.Cells.Item(10, 2).Value = 'Il16'
$ws.Cells.Item(10, 3).Value = 'Kcnj15'
$ws.Cells.Item(10, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(10, 5).Value = [double]"3"
$ws.Cells.Item(10, 6).Value = [double]"1"
$ws.Cells.Item(10, 7).Value = [double]"0.467591"
$ws.Cells.Item(10, 8).Value = [double]"1.402773"
$ws.Cells.Item(10, 9).Value = [double]"0.02856965900840602"
$ws.Cells.Item(10, 10).Value = [double]"0.02856965900840601"
$ws.Cells.Item(10, 11).Value = [double]"1"
$ws.Cells.Item(10, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(10, 13).Value = [double]"0.01106233333333333"
$ws.Cells.Item(10, 14).Value = [double]"0.033187"
$ws.Cells.Item(10, 15).Value = [double]"0.001140712431232264"
$ws.Cells.Item(10, 16).Value = [double]"0.001140712431232264"
$ws.Cells.Item(10, 17).Value = [double]"0.005172647505666667"
$ws.Cells.Item(10, 18).Value = [double]"0.046553827551"
$ws.Cells.Item(10, 19).Value = [double]"3.25897651869556E-05"
$ws.Cells.Item(10, 20).Value = [double]"3.258976518695559E-05"

# Row 11
$ws.Cells.Item(11, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(11, 2).Value = 'Il16'
$ws.Cells.Item(11, 3).Value = 'Kcnj15'
$ws.Cells.Item(11, 4).Value = 'ECs'
$ws.Cells.Item(11, 5).Value = [double]"3"
$ws.Cells.Item(11, 6).Value = [double]"1"
$ws.Cells.Item(11, 7).Value = [double]"5.892462666666667"
$ws.Cells.Item(11, 8).Value = [double]"17.677388"
$ws.Cells.Item(11, 9).Value = [double]"0.3600275649155554"
$ws.Cells.Item(11, 10).Value = [double]"0.3600275649155554"
$ws.Cells.Item(11, 11).Value = [double]"3"
$ws.Cells.Item(11, 12).Value = [double]"1"
$ws.Cells.Item(11, 13).Value = [double]"2.359186333333333"
$ws.Cells.Item(11, 14).Value = [double]"7.077559"
$ws.Cells.Item(11, 15).Value = [double]"0.243271749000506"
$ws.Cells.Item(11, 16).Value = [double]"0.243271749000506"
$ws.Cells.Item(11, 17).Value = [double]"13.90141739287689"
$ws.Cells.Item(11, 18).Value = [double]"125.112756535892"
$ws.Cells.Item(11, 19).Value = [double]"0.08758453540540038"
$ws.Cells.Item(11, 20).Value = [double]"0.08758453540540037"

# Row 12
$ws.Cells.Item(12, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(12, 2).Value = 'Il16'
$ws.Cells.Item(12, 3).Value = 'Kcnj15'
$ws.Cells.Item(12, 4).Value = 'FAPs'
$ws.Cells.Item(12, 5).Value = [double]"3"
$ws.Cells.Item(12, 6).Value = [double]"1"
$ws.Cells.Item(12, 7).Value = [double]"5.892462666666667"
$ws.Cells.Item(12, 8).Value = [double]"17.677388"
$ws.Cells.Item(12, 9).Value = [double]"0.3600275649155554"
$ws.Cells.Item(12, 10).Value = [double]"0.3600275649155554"
$ws.Cells.Item(12, 11).Value = [double]"3"
$ws.Cells.Item(12, 12).Value = [double]"1"
$ws.Cells.Item(12, 13).Value = [double]"7.327491999999999"
$ws.Cells.Item(12, 14).Value = [double]"21.982476"
$ws.Cells.Item(12, 15).Value = [double]"0.7555875385682617"
$ws.Cells.Item(12, 16).Value = [double]"0.7555875385682617"
$ws.Cells.Item(12, 17).Value = [double]"43.17697305029866"
$ws.Cells.Item(12, 18).Value = [double]"388.592757452688"
$ws.Cells.Item(12, 19).Value = [double]"0.2720323415912695"
$ws.Cells.Item(12, 20).Value = [double]"0.2720323415912695"

# Row 13
$ws.Cells.Item(13, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(13, 2).Value = 'Il16'
$ws.Cells.Item(13, 3).Value = 'Kcnj15'
$ws.Cells.Item(13, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(13, 5).Value = [double]"3"
$ws.Cells.Item(13, 6).Value = [double]"1"
$ws.Cells.Item(13, 7).Value = [double]"5.892462666666667"
$ws.Cells.Item(13, 8).Value = [double]"17.677388"
$ws.Cells.Item(13, 9).Value = [double]"0.3600275649155554"
$ws.Cells.Item(13, 10).Value = [double]"0.3600275649155554"
$ws.Cells.Item(13, 11).Value = [double]"1"
$ws.Cells.Item(13, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(13, 13).Value = [double]"0.01106233333333333"
$ws.Cells.Item(13, 14).Value = [double]"0.033187"
$ws.Cells.Item(13, 15).Value = [double]"0.001140712431232264"
$ws.Cells.Item(13, 16).Value = [double]"0.001140712431232264"
$ws.Cells.Item(13, 17).Value = [double]"0.06518438617288889"
$ws.Cells.Item(13, 18).Value = [double]"0.5866594755560001"
$ws.Cells.Item(13, 19).Value = [double]"0.0004106879188854552"
$ws.Cells.Item(13, 20).Value = [double]"0.0004106879188854551"
